$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '57.776.25'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').Value = '3.056.65'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.Value = '''524.81'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +5.64%  '
$c = $ws.Range('D6')
$c.Value = '''142.42'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +5.46%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +4.82%  '
$c = $ws.Range('D9')
$c.Value = '''7.64'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +5.10%  '
$ws.Range('E10').Value = '  +7.95%  '
$ws.Range('E11').Value = '  +5.04%  '
$ws.Range('E12').Value = '  +2.28%  '
$ws.Range('D13').Value = '3.577.19'
$ws.Range('E13').Value = '  +2.40%  '
$c = $ws.Range('D14')
$c.Value = '''26.89'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +7.92%  '
$c = $ws.Range('D15')
$c.Value = '''0.0000171'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +17.37%  '
$ws.Range('D16').Value = '57.754.85'
$ws.Range('E16').Value = '  +2.28%  '
$c = $ws.Range('D17')
$c.Value = '''6.23'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +6.85%  '
$ws.Range('D18').Value = '3.054.23'
$ws.Range('E18').Value = '  +2.44%  '
$c = $ws.Range('D19')
$c.Value = '''13.05'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +5.54%  '
$c = $ws.Range('D20')
$c.Value = '''8.19'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +5.56%  '
$c = $ws.Range('D21')
$c.Value = '''339.39'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +4.53%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  +7.48%  '
$c = $ws.Range('D24')
$c.Value = '''64.89'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +5.89%  '
$c = $ws.Range('D25')
$c.Value = '''0.172'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +5.89%  '
$ws.Range('D26').Value = '0.0₃0975'
$ws.Range('E26').Value = '  +9.23%  '
$ws.Range('E27').Value = '  +0.63%  '
$c = $ws.Range('D28')
$c.Value = '''6.91'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +6.45%  '
$c = $ws.Range('D29')
$c.Value = '''7.37'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +10.37%  '
$c = $ws.Range('D30')
$c.Value = '''1.86'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +7.36%  '
$ws.Range('E31').Value = '  +5.11%  '
$c = $ws.Range('D32')
$c.Value = '''21.12'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.69%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D33')
$c.Value = '''4.76'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +6.59%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D34')
$c.Value = '''156.50'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.59%  '
$c = $ws.Range('D35')
$c.Value = '''5.97'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +6.53%  '
$ws.Range('E36').Value = '  +3.24%  '
$c = $ws.Range('D37')
$c.Value = '''26.07'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +12.10%  '
$c = $ws.Range('D38')
$c.Value = '''0.0704'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.63%  '
$ws.Range('D39').Value = '3.092.62'
$ws.Range('E39').Value = '  +2.54%  '
$c = $ws.Range('D40')
$c.Value = '''37.73'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +3.24%  '
$ws.Range('E41').Value = '  +9.04%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  +5.69%  '
$c = $ws.Range('D44')
$c.Value = '''0.663'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +4.07%  '
$ws.Range('D45').Value = '2.330.84'
$ws.Range('E45').Value = '  +4.82%  '
$ws.Range('E46').Value = '  +4.41%  '
$c = $ws.Range('D47')
$c.Value = '''2.03'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +4.86%  '
$ws.Range('E48').Value = '  +4.78%  '
$ws.Range('E49').Value = '  +4.31%  '
$ws.Range('E50').Value = '  +6.21%  '
$ws.Range('E51').Value = '  +6.05%  '
